# Sari_Cases_11 export template update:
#  - Add five new symptom columns to the export sheet (with headers), inserted
#    among the existing "process*" / symptom columns rather than appended at
#    the end.
#  - Re-style the header row: bold, 12pt, vertically centered, taller row.
#  - Leave active selection on A2 (matches the saved workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Asymptomatic" header column ---------------------------------
$ws.Range("CO1").EntireColumn.Insert()
$ws.Range("CO1").Value = "Asymptomatic"

# --- Insert four more new symptom header columns --------------------------
$ws.Range("CZ1:DC1").EntireColumn.Insert()
$ws.Range("CZ1").Value = "Rhinorrhoea "
$ws.Range("DA1").Value = "Odinophagia"
$ws.Range("DB1").Value = "Anosmy"
$ws.Range("DC1").Value = "Dysgeusia"

# --- Reformat the header row ----------------------------------------------
$headerRow = $ws.Rows.Item(1)
$headerRow.Font.Bold = $true
$headerRow.Font.Size = 12
$headerRow.VerticalAlignment = -4108
$headerRow.RowHeight = 25.5

# --- Restore the saved selection state ------------------------------------
$ws.Range("A2").Select()
